$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (week 14->B15, 22->B23, 23->B24, 24->B25, 25->B26)
$ws.Range("B15").Value = 457
$ws.Range("B23").Value = 347
$ws.Range("B24").Value = 359
$ws.Range("B25").Value = 263
$ws.Range("B26").Value = 353

# Correct previously placeholder value for week 35 and add week 36
$ws.Range("B36").Value = 412

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 439
